# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (fund-holding detail) right before the
#   "总计" (totals) summary sheet, cloned from the "2021-Q4" sheet so it
#   inherits the same header/column formatting, then filled with the new
#   quarter's data.
# - Insert a new top data row into "总计" for "2022-Q1" (6 funds held,
#   0.92 亿元 total market value) and renumber the existing index column.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q1" sheet by cloning "2021-Q4" (same layout:
#        bold/bordered header row + bold/bordered index column A) and
#        dropping it in right before "总计". Worksheet handles in this host
#        are position-anchored, so any handle obtained before a sheet-level
#        structural change (Copy/Add/Move/Delete) can silently end up
#        pointing at a different sheet afterwards -- always re-fetch by
#        name right after such an operation. ---
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$template.Copy($totalSheet)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Extend the cloned formatting from 5 data rows to 7 data rows (rows 6-7),
# reusing row 5's per-column styles (bold+border index cell, plain data
# cells) so no stray styles get created.
$q1.Range("A5:H5").Copy()
$q1.Range("A6:H7").PasteSpecial(-4122)
$q1.Range("A1").Select()

# --- 2. Fill in the 2022-Q1 holdings data. Numeric-looking text columns
#        (fund size / position / weight / market value) are entered with a
#        leading apostrophe so they stay text, matching the workbook's
#        existing convention (only the rank column H and index column A
#        are real numbers). ---
$q1Data = @(
    @(0, "513060", "博时恒生医疗保健ETF（QDII）", "23.35", "98.89", "2.77", "0.6468", 9),
    @(1, "159792", "富国中证港股通互联网ETF", "2.76", "99.00", "3.58", "0.0988", 8),
    @(2, "513700", "鹏华中证港股通医药卫生综合交易型开放式指数证券投资基金", "3.24", "93.11", "2.43", "0.0787", 9),
    @(3, "003993", "前海开源沪港深核心驱动灵活配置混合", "0.58", "82.10", "7.98", "0.0463", 5),
    @(4, "159892", "华夏恒生香港上市生物科技ETF（QDII）", "1.51", "99.03", "3.02", "0.0456", 8),
    @(5, "004321", "前海开源沪港深强国产业灵活配置混合", "0.02", "64.32", "4.14", "0.0008", 7)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- 3. Insert the new 2022-Q1 summary row at the top of "总计" and
#        renumber the existing rows' index column. Re-fetch the "总计"
#        handle by name since the sheet-copy above invalidated the old one. ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.92

for ($row = 3; $row -le 7; $row++) {
    $totalSheet.Cells.Item($row, 1).Value = $row - 2
}

$totalSheet.Range("A1").Select()
